$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking values are not
# converted to Excel numbers (the source data stores prices as text).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '66.858.06'
$ws.Range("E2").Value = '  +0.80%  '

$ws.Range("D3").Value = '3.501.43'
$ws.Range("E3").Value = '  +0.39%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = '595.07'
$ws.Range("E5").Value = '  +0.48%  '

$ws.Range("D6").Value = '169.88'
$ws.Range("E6").Value = '  +0.40%  '

$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("D8").Value = '0.591'
$ws.Range("E8").Value = '  +2.54%  '

$ws.Range("D9").Value = '0.134'
$ws.Range("E9").Value = '  +7.62%  '

$ws.Range("E10").Value = '  +0.81%  '

$ws.Range("E11").Value = '  -0.51%  '

$ws.Range("D12").Value = '4.103.30'
$ws.Range("E12").Value = '  +0.31%  '

$ws.Range("E13").Value = '  -0.18%  '

$ws.Range("D14").Value = '28.31'
$ws.Range("E14").Value = '  +1.53%  '

$ws.Range("E15").Value = '  +3.28%  '

$ws.Range("D16").Value = '66.809.88'
$ws.Range("E16").Value = '  +0.77%  '

$ws.Range("D17").Value = '3.488.16'
$ws.Range("E17").Value = '  -0.23%  '

$ws.Range("D18").Value = '6.33'
$ws.Range("E18").Value = '  +0.93%  '

$ws.Range("D19").Value = '14.10'
$ws.Range("E19").Value = '  +0.62%  '

$ws.Range("D20").Value = '395.53'
$ws.Range("E20").Value = '  +2.10%  '

$ws.Range("D21").Value = '7.96'
$ws.Range("E21").Value = '  -0.27%  '

$ws.Range("D22").Value = '73.35'
$ws.Range("E22").Value = '  +0.42%  '

$ws.Range("E23").Value = '  -0.16%  '

$ws.Range("D24").Value = '0.537'
$ws.Range("E24").Value = '  +2.04%  '

$ws.Range("E25").Value = '  +1.20%  '

$ws.Range("D26").Value = '10.22'
$ws.Range("E26").Value = '  +1.38%  '

$ws.Range("E27").Value = '  +0.21%  '

$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  -0.01%  '

$ws.Range("D29").Value = '6.30'
$ws.Range("E29").Value = '  -1.12%  '

$ws.Range("D30").Value = '1.46'
$ws.Range("E30").Value = '  -0.89%  '

$ws.Range("E31").Value = '  +0.23%  '

$ws.Range("D32").Value = '24.10'
$ws.Range("E32").Value = '  +2.93%  '

$ws.Range("D33").Value = '7.38'
$ws.Range("E33").Value = '  -0.17%  '

$ws.Range("E34").Value = '  +4.53%  '

$ws.Range("D35").Value = '163.70'
$ws.Range("E35").Value = '  +1.81%  '

$ws.Range("D36").Value = '0.893'
$ws.Range("E36").Value = '  -0.84%  '

$ws.Range("B38").Value = 'Filecoin'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D38").Value = '4.73'
$ws.Range("E38").Value = '  +3.42%  '

$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = '6.85'
$ws.Range("E39").Value = '  +3.03%  '

$ws.Range("D40").Value = '0.0744'
$ws.Range("E40").Value = '  -0.23%  '

$ws.Range("D41").Value = '26.45'
$ws.Range("E41").Value = '  +0.33%  '

$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").Value = '2.64'
$ws.Range("E42").Value = '  +6.48%  '

$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '2.824.53'
$ws.Range("E43").Value = '  +1.23%  '

$ws.Range("D44").Value = '26.96'
$ws.Range("E44").Value = '  -0.81%  '

$ws.Range("D45").Value = '42.79'
$ws.Range("E45").Value = '  -1.37%  '

$ws.Range("E46").Value = '  -0.46%  '

$ws.Range("D47").Value = '341.66'
$ws.Range("E47").Value = '  -3.11%  '

$ws.Range("E48").Value = '  +1.81%  '

$ws.Range("D49").Value = '33.48'
$ws.Range("E49").Value = '  +2.39%  '

$ws.Range("E50").Value = '  +0.37%  '

$ws.Range("E51").Value = '  +0.97%  '

# Restore default styling on column D so no stray number-format style is left
# behind on cells (matches original workbook which has no explicit style here).
$ws.Range("D2:D51").Style = "Normal"
